$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update names (column A) and counts (column B) for rows 2-8 and 11
# to reflect the new sorted order and refreshed totals.

$ws.Range("A2").Value = "PALACIOS PANTA LUIS MIGUEL"
$ws.Range("B2").Value = 109

$ws.Range("A3").Value = "VEGA ZAPATA JESUS GABRIEL"
$ws.Range("B3").Value = 106

$ws.Range("A4").Value = "PANTA NIMA FREDDY ROLAND JUNIOR"
$ws.Range("B4").Value = 105

$ws.Range("A5").Value = "CRISANTO CARMEN ROSITA ABIGAIL"
$ws.Range("B5").Value = 104

$ws.Range("A6").Value = "SALAZAR VEGA MARIA FERNANDA"
$ws.Range("B6").Value = 103

$ws.Range("A7").Value = "MAZA RIOFRIO CINTHIA NATELAHI"
$ws.Range("B7").Value = 100

$ws.Range("A8").Value = "PANTA VARONA CANDY ELIZABETH"
$ws.Range("B8").Value = 98

$ws.Range("B11").Value = 82
